$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per repulled data
$updates = @{
    4  = -1
    7  = 1
    14 = -2
    22 = -6
    26 = -6
    27 = -3
    28 = -2
    35 = -4
    38 = -3
    40 = 2
    43 = 2
    46 = -7
    48 = 1
    53 = -2
    58 = 0
    61 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
